$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D as text so numeric-looking price strings (e.g. "161.80",
# "1.00") keep their exact characters instead of being coerced to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '68.328.13'
$ws.Range('E2').Value = '  -2.01%  '
$ws.Range('D3').Value = '2.438.47'
$ws.Range('E3').Value = '  -2.85%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '560.82'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').Value = '161.80'
$ws.Range('E6').Value = '  -3.21%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -2.69%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '2.437.46'
$ws.Range('E9').Value = '  -2.91%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  -6.63%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.164'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').Value = '0.339'
$ws.Range('E12').Value = '  -6.32%  '
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D13').Value = '4.78'
$ws.Range('E13').Value = '  -3.27%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.885.22'
$ws.Range('E14').Value = '  -3.20%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '68.171.24'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000170'
$ws.Range('E16').Value = '  -4.43%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '23.51'
$ws.Range('E17').Value = '  -5.53%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.447.16'
$ws.Range('E18').Value = '  -2.85%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '10.91'
$ws.Range('E19').Value = '  -3.62%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '348.01'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '7.17'
$ws.Range('E21').Value = '  -5.60%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = '3.77'
$ws.Range('E22').Value = '  -3.79%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '1.83'
$ws.Range('E24').Value = '  -7.42%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '67.49'
$ws.Range('E25').Value = '  -5.19%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = '3.72'
$ws.Range('E26').Value = '  -5.91%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.564.82'
$ws.Range('E27').Value = '  -3.32%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '8.21'
$ws.Range('E29').Value = '  -7.35%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0831'
$ws.Range('E30').Value = '  -7.01%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '7.29'
$ws.Range('E31').Value = '  -7.67%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -5.32%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = '427.03'
$ws.Range('E34').Value = '  -7.48%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').Value = '1.67'
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('B36').Value = 'POPCAT'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range('D36').Value = '3.17'
$ws.Range('E36').Value = '  +110.87%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '157.42'
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').Value = '18.98'
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.109'
$ws.Range('E40').Value = '  -6.19%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').Value = '17.88'
$ws.Range('E41').Value = '  -3.63%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = '0.304'
$ws.Range('E42').Value = '  -4.74%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').Value = '4.46'
$ws.Range('E43').Value = '  -5.05%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.52'
$ws.Range('E44').Value = '  -5.25%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').Value = '1.06'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = '2.05'
$ws.Range('E46').Value = '  -7.48%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '134.97'
$ws.Range('E47').Value = '  -5.42%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').Value = '3.33'
$ws.Range('E48').Value = '  -4.22%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0715'
$ws.Range('E49').Value = '  -2.63%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '0.483'
$ws.Range('E50').Value = '  -7.53%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.556'
$ws.Range('E51').Value = '  -3.89%  '
